# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values on the zh-cn and de-de
# sheets to reflect the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 3 (9c330023-... file)
$wsZhCn.Range("D3").Value = "2016-02-26 06:33:12"
$wsZhCn.Range("G3").Value = "2016-02-26 06:34:08"

# de-de sheet, row 2 (57ebf84e-... file)
$wsDeDe.Range("G2").Value = "2016-02-26 06:33:27"

# de-de sheet, row 3 (9c330023-... file)
$wsDeDe.Range("D3").Value = "2016-02-26 06:33:27"
$wsDeDe.Range("G3").Value = "2016-02-26 06:34:32"

Write-Host "Handback report timestamps updated."
